# Updated cryptos list - refresh Price (D) and Volume(1h) (E) columns,
# plus the USDC/BNB row swap (rows 5 and 6), per upstream coinranking data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values are written with a leading apostrophe so Excel keeps them as
# plain text (these columns hold text like "26.486.49" / "0.9978" / "  -3.15%  ",
# not numbers), then the style is reset back to Normal so no Text number-
# format is left behind on the cell.
function Set-TextValue($ref, $val) {
    $ws.Range($ref).Value = "'" + $val
    $ws.Range($ref).Style = "Normal"
}

Set-TextValue "D2" "26.486.49"
Set-TextValue "E2" "  -3.15%  "
Set-TextValue "D3" "1.774.90"
Set-TextValue "E3" "  -2.09%  "
Set-TextValue "D4" "0.9984"
Set-TextValue "E4" "  -0.45%  "
Set-TextValue "B5" "USDC"
Set-TextValue "C5" "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
Set-TextValue "D5" "0.9978"
Set-TextValue "E5" "  -0.44%  "
Set-TextValue "B6" "BNB"
Set-TextValue "C6" "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
Set-TextValue "D6" "304.77"
Set-TextValue "E6" "  -2.14%  "
Set-TextValue "D7" "0.4280"
Set-TextValue "E7" "  +1.43%  "
Set-TextValue "D8" "0.3642"
Set-TextValue "E8" "  +2.02%  "
Set-TextValue "D9" "0.07182"
Set-TextValue "E9" "  +0.29%  "
Set-TextValue "D10" "0.8493"
Set-TextValue "E10" "  -0.07%  "
Set-TextValue "D11" "20.49"
Set-TextValue "E11" "  +1.06%  "
Set-TextValue "D12" "1.771.83"
Set-TextValue "E12" "  -5.17%  "
Set-TextValue "D13" "6.457"
Set-TextValue "E13" "  +1.01%  "
Set-TextValue "D14" "5.264"
Set-TextValue "E14" "  -1.30%  "
Set-TextValue "D15" "0.06866"
Set-TextValue "E15" "  -0.74%  "
Set-TextValue "D16" "0.9990"
Set-TextValue "E16" "  -0.63%  "
Set-TextValue "D17" "78.96"
Set-TextValue "E17" "  -3.31%  "
Set-TextValue "D18" "0.000008698"
Set-TextValue "E18" "  -1.70%  "
Set-TextValue "D19" "0.9979"
Set-TextValue "E19" "  -0.47%  "
Set-TextValue "D20" "15.02"
Set-TextValue "E20" "  -1.14%  "
Set-TextValue "D21" "26.495.19"
Set-TextValue "E21" "  -3.87%  "
Set-TextValue "D22" "5.115"
Set-TextValue "E22" "  +0.22%  "
Set-TextValue "D23" "11.13"
Set-TextValue "E23" "  +1.34%  "
Set-TextValue "D24" "2.012.84"
Set-TextValue "E24" "  -3.44%  "
Set-TextValue "D25" "152.33"
Set-TextValue "E25" "  -1.04%  "
Set-TextValue "D26" "1.870"
Set-TextValue "E26" "  -5.09%  "
Set-TextValue "E27" "  -1.11%  "
Set-TextValue "D28" "5.090"
Set-TextValue "E28" "  -0.30%  "
Set-TextValue "D29" "113.91"
Set-TextValue "E29" "  +0.41%  "
Set-TextValue "D30" "1.807"
Set-TextValue "E30" "  +4.12%  "
Set-TextValue "D31" "0.08933"
Set-TextValue "E31" "  +0.26%  "
Set-TextValue "D32" "0.7292"
Set-TextValue "E32" "  -1.95%  "
Set-TextValue "D33" "1.131"
Set-TextValue "E33" "  +1.39%  "
Set-TextValue "D34" "4.334"
Set-TextValue "E34" "  -3.41%  "
Set-TextValue "E35" "  -6.30%  "
Set-TextValue "D36" "0.9984"
Set-TextValue "E36" "  -0.39%  "
Set-TextValue "D37" "1.104"
Set-TextValue "E37" "  +2.80%  "
Set-TextValue "D38" "0.05160"
Set-TextValue "E38" "  -1.11%  "
Set-TextValue "D39" "0.01897"
Set-TextValue "E39" "  -0.74%  "
Set-TextValue "D40" "0.4941"
Set-TextValue "E40" "  -1.24%  "
Set-TextValue "D41" "0.1614"
Set-TextValue "E41" "  -1.92%  "
Set-TextValue "D42" "2.618"
Set-TextValue "E42" "  -5.47%  "
Set-TextValue "D43" "6.323"
Set-TextValue "E43" "  +0.00%  "
Set-TextValue "D44" "8.043"
Set-TextValue "E44" "  -2.58%  "
Set-TextValue "D45" "105.10"
Set-TextValue "E45" "  -0.17%  "
Set-TextValue "D46" "10.13"
Set-TextValue "E46" "  -1.82%  "
Set-TextValue "D47" "0.9969"
Set-TextValue "E47" "  -0.51%  "
Set-TextValue "D48" "1.633"
Set-TextValue "E48" "  +1.65%  "
Set-TextValue "D49" "0.4500"
Set-TextValue "E49" "  -2.39%  "
Set-TextValue "D50" "0.06202"
Set-TextValue "E50" "  -3.37%  "
Set-TextValue "D51" "1.738"
Set-TextValue "E51" "  +2.74%  "
